$d = $word.ActiveDocument

# Step 1: Replace the text of the "Calcular factor de concentración..." paragraph.
# The new wording is split across two sentences/runs in the final document, so
# we build it as two separate insertions to mirror the two <w:r> runs.
$oldText = "Calcular factor de concentración de tensiones, teóricamente y con los datos experimentales"
$firstPart = "Cuanta fuerza se necesita para realizar la carga máxima encontrada teóricamente"
$secondPart = ", para poder comparar y analizar los resultados"

$matchRange = $d.Content
$found = $matchRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($matchRange.Find.Found) {
    $matchRange.Text = $firstPart
    $insertPoint = $d.Range($matchRange.End, $matchRange.End)
    $insertPoint.InsertAfter($secondPart)
}

# Step 2: Remove the three paragraphs that followed entirely (they were
# deleted wholesale in the diff): "Obtener valor máximo para la tensión
# analítica", "Obtener valor máximo para la tensión con los datos
# experimentales" and "Buscar ecuaciones de Lamé". Deleting the paragraph's
# own Range (including its paragraph mark) removes the whole block.
$targets = @(
    "Obtener valor máximo para la tensión analítica",
    "Obtener valor máximo para la tensión con los datos experimentales",
    "Buscar ecuaciones de Lamé"
)

foreach ($target in $targets) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        $text = $para.Range.Text
        if ($text.TrimEnd([char]13, [char]7) -eq $target) {
            $para.Range.Delete()
            break
        }
    }
}
